$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell H1 onto the new header
# cells (same bold/centered/bordered header style used by the rest of row 1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-13 for columns I (I0) and J (IF)
$data = @{
    2  = @(1, 6)
    3  = @(1, 4)
    4  = @(1, 6)
    5  = @(4, 8)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(2, 5)
    9  = @(5, 9)
    10 = @(6, 9)
    11 = @(5, 6)
    12 = @(5, 6)
    13 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
